$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns: A=NUMBER, B=DATE, C=TIME, D=FLIGHT, E=FROM, F=SHORT,
#          G=AIRLINE, H=MODEL, I=AIRCFAT ID, J=STATUS, K=(blank),
#          L=DIFFERENCE, M=(blank)
$data = @(
    @(290, "Sunday, Jan 15", "2:50 AM",  "W61084",  "Fuerteventura", "(FUE)", "Wizz Air ",       "A21N", "(HA-LVT)", "2:23 AM",  "0 hours, -27 minutes"),
    @(291, "Sunday, Jan 15", "3:48 AM",  "BO624",   "Catania",       "(CTA)", "Bluebird Nordic ", "B734", "(TF-BBN)", "3:36 AM",  "0 hours, -12 minutes"),
    @(292, "Sunday, Jan 15", "6:10 AM",  "3Z7109",  "Marsa Alam",    "(RMF)", "Smartwings ",      "B38M", "(OK-SWC)", "5:55 AM",  "0 hours, -15 minutes"),
    @(293, "Sunday, Jan 15", "8:25 AM",  "FR3408",  "Milan",         "(BGY)", "Ryanair ",         "B738", "(9H-QAK)", "8:08 AM",  "0 hours, -17 minutes"),
    @(294, "Sunday, Jan 15", "8:34 AM",  "BO950",   "Cagliari",      "(CAG)", "Bluebird Nordic ", "B734", "(TF-BBJ)", "8:11 AM",  "0 hours, -23 minutes"),
    @(295, "Sunday, Jan 15", "9:25 AM",  "FR6392",  "London",        "(STN)", "Ryanair ",         "B38M", "(EI-HMS)", "9:17 AM",  "0 hours, -8 minutes"),
    @(296, "Sunday, Jan 15", "11:10 AM", "LO3881",  "Warsaw",        "(WAW)", "LOT ",             "E170", "(SP-LDI)", "10:59 AM", "0 hours, -11 minutes"),
    @(297, "Sunday, Jan 15", "11:30 AM", "FR6372",  "Catania",       "(CTA)", "Ryanair ",         "B738", "(SP-RSH)", "11:10 AM", "0 hours, -20 minutes"),
    @(298, "Sunday, Jan 15", "11:40 AM", "FR6893",  "Dortmund",      "(DTM)", "Ryanair ",         "B738", "(SP-RSB)", "11:29 AM", "0 hours, -11 minutes"),
    @(299, "Sunday, Jan 15", "11:50 AM", "W61002",  "London",        "(LTN)", "Wizz Air ",        "A21N", "(HA-LZF)", "11:38 AM", "0 hours, -12 minutes"),
    @(300, "Sunday, Jan 15", "12:00 PM", "W61020",  "Bristol",       "(BRS)", "Wizz Air ",        "A21N", "(HA-LZJ)", "11:51 AM", "0 hours, -9 minutes"),
    @(301, "Sunday, Jan 15", "12:25 PM", "W61242",  "Athens",        "(ATH)", "Wizz Air ",        "A21N", "(HA-LVT)", "12:09 PM", "0 hours, -16 minutes"),
    @(302, "Sunday, Jan 15", "1:20 PM",  "FR5397",  "Dublin",        "(DUB)", "Ryanair ",         "B738", "(EI-EVC)", "1:08 PM",  "0 hours, -12 minutes"),
    @(303, "Sunday, Jan 15", "1:33 PM",  "UNKNOWN", "Cologne",       "(CGN)", "Enter Air ",       "B738", "(SP-ENW)", "3:16 PM",  "1 hours, 43 minutes"),
    @(304, "Sunday, Jan 15", "1:40 PM",  "LH1356",  "Frankfurt",     "(FRA)", "Lufthansa ",       "CRJ9", "(D-ACNB)", "1:26 PM",  "0 hours, -14 minutes"),
    @(305, "Sunday, Jan 15", "2:40 PM",  "LO3883",  "Warsaw",        "(WAW)", "LOT ",             "E75S", "(SP-LID)", "2:27 PM",  "0 hours, -13 minutes"),
    @(306, "Sunday, Jan 15", "2:47 PM",  "E44832",  "Hurghada",      "(HRG)", "Enter Air ",       "B738", "(SP-ENP)", "2:14 PM",  "0 hours, -33 minutes"),
    @(307, "Sunday, Jan 15", "4:10 PM",  "KL1815",  "Amsterdam",     "(AMS)", "KLM ",             "E295", "(PH-NXN)", "3:54 PM",  "0 hours, -16 minutes"),
    @(308, "Sunday, Jan 15", "4:25 PM",  "3Z7111",  "Hurghada",      "(HRG)", "Smartwings ",      "B738", "(OK-TSF)", "3:12 PM",  "-1 hours, -13 minutes")
)

$startRow = 291
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value2 = [double]$rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $ws.Cells.Item($row, 10).Value = $rec[9]
    # Columns K and M stay blank in the source data, but the row still
    # carries a (styled) cell there, so touch a formatting property to
    # materialize the cell without introducing a new style.
    $ws.Cells.Item($row, 11).Font.Size = 11
    $ws.Cells.Item($row, 12).Value = $rec[10]
    $ws.Cells.Item($row, 13).Font.Size = 11
}
